# EPBDS-8844 Customizing output of a SpreadsheetResult.
#
# Renames the asterisk-prefixed section markers ("*Step1", "*Values",
# "*Formula", "*Step2", "*moreValues") to asterisk-suffixed markers
# ("Step1*", "Values*", "Formula*", "Step2*", "moreValues*"), and renames
# the "Ignored" spreadsheet samples (sprTwoTwoIgnored1 / sprTwoTwoIgnored2)
# to sprTwoTwo1 / sprTwoTwo2. Also widens column C and refreshes the
# worksheet selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-Literal($addr, $text) {
    # Leading apostrophe forces Excel to store the value as literal text
    # instead of (re)parsing a leading "=" as a formula.
    $ws.Range($addr).Value = "'" + $text
}

# ---- sprOneRow / "B5" block header & sprOneColumn block (rows 5-9) ----
Set-Literal "B5" "Spreadsheet SpreadsheetResult sprOneColumn(Integer myVar)"

Set-Literal "B6" "properties"
Set-Literal "C6" "state"
Set-Literal "D6" "AL"

Set-Literal "C7" "Step1*"

Set-Literal "B8" "Values*"
Set-Literal "C8" "=myVar + 1"
Set-Literal "D8" "=myVar + 3"

Set-Literal "B9" "moreValues*"
Set-Literal "C9" "=myVar + 2"
Set-Literal "D9" "=myVar + 4"

# ---- rows 20-25 ----
Set-Literal "B20" "Spreadsheet SpreadsheetResult sprOneRow(Integer myVar)"

Set-Literal "B21" "properties"
Set-Literal "C21" "state"
Set-Literal "D21" "AL"

Set-Literal "C22" "Values*"
Set-Literal "D22" "moreValues*"

Set-Literal "B23" "Step1*"
Set-Literal "C23" "=myVar + 1"
Set-Literal "D23" "=myVar + 3"

Set-Literal "C24" "=myVar + 2"
Set-Literal "D24" "=myVar + 4"

# ---- sprTwoTwo block (rows 32-36) ----
Set-Literal "B32" "Spreadsheet SpreadsheetResult sprTwoTwo(Integer myVar)"

Set-Literal "B33" "properties"
Set-Literal "C33" "state"
Set-Literal "D33" "AL"

Set-Literal "C34" "Values*"
Set-Literal "D34" "Formula*"

Set-Literal "B35" "Step1*"
Set-Literal "C35" "=myVar + 1"
Set-Literal "D35" "=myVar + 3"

Set-Literal "B36" "Step2*"
Set-Literal "C36" "=myVar + 2"
Set-Literal "D36" "=myVar + 4"

# ---- sprOneOne block (rows 43-46) ----
Set-Literal "B43" "Spreadsheet SpreadsheetResult sprOneOne(Integer myVar)"

Set-Literal "C44" "Values*"

Set-Literal "B45" "Step1*"
Set-Literal "C45" "=myVar + 1"

# ---- dtRetSpr block (rows 48-50) ----
Set-Literal "B48" "SimpleRules SpreadsheetResult dtRetSpr(Integer v)"

Set-Literal "B49" "Value"
Set-Literal "C49" "Ret"

Set-Literal "C50" "=sprOneOne(v)"

# ---- dtRetSpr2 block (rows 53-55) ----
Set-Literal "B53" "SimpleRules SpreadsheetResult[] dtRetSpr2(Integer[] v)"

Set-Literal "B54" "Value"
Set-Literal "C54" "Ret"

Set-Literal "C55" "=sprOneOne(v)"

# ---- sprOneOneNoAsterisk block (rows 58-60) ----
Set-Literal "B58" "Spreadsheet SpreadsheetResult sprOneOneNoAsterisk(Integer myVar)"

Set-Literal "C59" "Values"

Set-Literal "B60" "Step1"
Set-Literal "C60" "=myVar + 1"

# ---- sprTwoTwoIgnored1 -> sprTwoTwo1 block (rows 63-67) ----
Set-Literal "B63" "Spreadsheet SpreadsheetResult sprTwoTwo1(Integer myVar)"

Set-Literal "B64" "properties"
Set-Literal "C64" "state"
Set-Literal "D64" "AL"

Set-Literal "D65" "Formula"

Set-Literal "B66" "Step1*"
Set-Literal "C66" "=myVar + 1"
Set-Literal "D66" "=myVar + 3"

Set-Literal "B67" "Step2"
Set-Literal "C67" "=myVar + 2"
Set-Literal "D67" "=myVar + 4"

# ---- sprTwoTwoIgnored2 -> sprTwoTwo2 block (rows 70-74) ----
Set-Literal "B70" "Spreadsheet SpreadsheetResult sprTwoTwo2(Integer myVar)"

Set-Literal "B71" "properties"
Set-Literal "C71" "state"
Set-Literal "D71" "AL"

Set-Literal "C72" "Values*"
Set-Literal "D72" "Formula*"

Set-Literal "B73" "Step1"
Set-Literal "C73" "=myVar + 1"
Set-Literal "D73" "=myVar + 3"

Set-Literal "C74" "=myVar + 2"
Set-Literal "D74" "=myVar + 4"

# ---- dtRetSpr3 block (rows 78-80) ----
Set-Literal "B78" "SimpleRules SpreadsheetResultsprOneOne[] dtRetSpr3(Integer[] v)"

Set-Literal "B79" "Value"
Set-Literal "C79" "Ret"

Set-Literal "C80" "=sprOneOne(v)"

# ---- Column C widened to fit the longer "*"-suffixed labels ----
$ws.Columns.Item(3).ColumnWidth = 43

# ---- Refresh view: drop the old scrolled/selected state, select C28 ----
$ws.Range("C28").Select()
